# fixed study time calculation
# Rolls the 30-day "Monthly Data" window forward to end on 2024-08-14
# (previously it ended on 2024-08-01), and corrects the Hours Studied
# value that was recorded for 2024-08-01 (3 -> 2). Also refreshes the
# "Today Progress" remaining-hours figure to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly Data")

# New 30-row window: row 2 = most recent day (2024-08-14) down to
# row 31 = oldest day (2024-07-16), mirroring the date-descending layout
# already used by the sheet.
$rows = @(
    @{ Date = "2024-08-14"; Hours = 0.5 },
    @{ Date = "2024-08-13"; Hours = 0 },
    @{ Date = "2024-08-12"; Hours = 0 },
    @{ Date = "2024-08-11"; Hours = 0 },
    @{ Date = "2024-08-10"; Hours = 0 },
    @{ Date = "2024-08-09"; Hours = 0 },
    @{ Date = "2024-08-08"; Hours = 0 },
    @{ Date = "2024-08-07"; Hours = 0 },
    @{ Date = "2024-08-06"; Hours = 0 },
    @{ Date = "2024-08-05"; Hours = 2.833333333333333 },
    @{ Date = "2024-08-04"; Hours = 0 },
    @{ Date = "2024-08-03"; Hours = 0 },
    @{ Date = "2024-08-02"; Hours = 3 },
    @{ Date = "2024-08-01"; Hours = 2 },
    @{ Date = "2024-07-31"; Hours = 2.25 },
    @{ Date = "2024-07-30"; Hours = 2 },
    @{ Date = "2024-07-29"; Hours = 2.083333333333333 },
    @{ Date = "2024-07-28"; Hours = 0.75 },
    @{ Date = "2024-07-27"; Hours = 3 },
    @{ Date = "2024-07-26"; Hours = 0 },
    @{ Date = "2024-07-25"; Hours = 1.75 },
    @{ Date = "2024-07-24"; Hours = 4 },
    @{ Date = "2024-07-23"; Hours = 0 },
    @{ Date = "2024-07-22"; Hours = 1.5 },
    @{ Date = "2024-07-21"; Hours = 3 },
    @{ Date = "2024-07-20"; Hours = 0 },
    @{ Date = "2024-07-19"; Hours = 0 },
    @{ Date = "2024-07-18"; Hours = 0 },
    @{ Date = "2024-07-17"; Hours = 0 },
    @{ Date = "2024-07-16"; Hours = 0 }
)

# Force column A to stay plain text (it was authored as literal date
# strings, not real date serials) instead of letting Excel auto-convert
# the ISO-looking strings into date values.
$dateRange = $ws.Range("A2:A31")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i].Date
    $ws.Cells.Item($r, 2).Value = $rows[$i].Hours
}

# Restore the default (unstyled) look now that the text is committed.
$dateRange.Style = "Normal"

# "Today Progress" sheet: remaining hours for today's goal changes
# from 3 to 0.5 (completed hours stays 0).
$progress = $wb.Worksheets.Item("Today Progress")
$progress.Range("B3").Value = 0.5

Write-Output "study window rolled forward; today remaining hours updated"
